# Fruta / hortaliza, semanal
# New weekly price observation inserted at row 133 ("Poroto verde" / Magnum /
# Primera / Provincia de Limarí), pushing every subsequent record down one
# row (old row 133 -> 134, ..., old row 219 -> 220).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 133; Excel shifts rows 133-219
# down to 134-220 and copies the surrounding row formatting (incl. the date
# style on column D) automatically.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(133, 1).Value = 8
$ws.Cells.Item(133, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(133, 3).Value = "Coquimbo"
$ws.Cells.Item(133, 4).Value = 44680
$ws.Cells.Item(133, 5).Value = 4
$ws.Cells.Item(133, 6).Value = 100112031
$ws.Cells.Item(133, 7).Value = "Poroto verde"
$ws.Cells.Item(133, 8).Value = "Magnum"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 400
$ws.Cells.Item(133, 11).Value = 23000
$ws.Cells.Item(133, 12).Value = 24000
$ws.Cells.Item(133, 13).Value = 23500
$ws.Cells.Item(133, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(133, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(133, 16).Value = 940
$ws.Cells.Item(133, 17).Value = 25
$ws.Cells.Item(133, 18).Value = "Hortaliza"
